# Auto-update predictions and index for 2025-10-12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Force the cell to stay text-typed (matches the workbook's existing
    # convention of storing every data value, even numeric-looking ones,
    # as inline/shared text) instead of Excel auto-coercing to a number.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- New "Odds" header column ---
Set-TextValue $ws.Range("G1") "Odds"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 2: Croatia vs Gibraltar ---
Set-TextValue $ws.Range("E2") "46 of 47 Tips"
Set-TextValue $ws.Range("G2") "1.01"

# --- Row 3: Faroe Islands vs Czech Republic ---
Set-TextValue $ws.Range("G3") "1.42"

# --- Row 4: now San Marino vs Cyprus ---
Set-TextValue $ws.Range("A4") "San Marino vs Cyprus"
Set-TextValue $ws.Range("B4") "Cyprus"
Set-TextValue $ws.Range("D4") "12th Oct 14:00"
Set-TextValue $ws.Range("E4") "20 of 21 Tips"
Set-TextValue $ws.Range("F4") "95"
Set-TextValue $ws.Range("G4") "1.16"

# --- Row 5: now Denmark vs Greece ---
Set-TextValue $ws.Range("A5") "Denmark vs Greece"
Set-TextValue $ws.Range("B5") "Denmark"
Set-TextValue $ws.Range("E5") "19 of 28 Tips"
Set-TextValue $ws.Range("F5") "68"
Set-TextValue $ws.Range("G5") "1.70"

# --- Row 6: now Lithuania vs Poland ---
Set-TextValue $ws.Range("A6") "Lithuania vs Poland"
Set-TextValue $ws.Range("B6") "Poland"
Set-TextValue $ws.Range("D6") "12th Oct 19:45"
Set-TextValue $ws.Range("E6") "19 of 22 Tips"
Set-TextValue $ws.Range("F6") "86"
Set-TextValue $ws.Range("G6") "1.44"

# --- Row 7: Romania vs Austria ---
Set-TextValue $ws.Range("G7") "1.84"

# --- Row 8: Netherlands vs Finland ---
Set-TextValue $ws.Range("G8") "1.11"

# --- Row 9: Universidad de Chile vs Palestino ---
Set-TextValue $ws.Range("G9") "1.74"

# --- Row 10: Nottingham Forest vs Chelsea ---
Set-TextValue $ws.Range("G10") "1.96"

# --- Row 11: Scotland vs Belarus ---
Set-TextValue $ws.Range("G11") "1.18"

# --- Row 12: Fulham vs Arsenal ---
Set-TextValue $ws.Range("G12") "1.60"

# --- Row 13: Man City vs Everton ---
Set-TextValue $ws.Range("G13") "1.42"

# --- Row 14: now Barcelona vs Olympiacos ---
Set-TextValue $ws.Range("A14") "Barcelona vs Olympiacos"
Set-TextValue $ws.Range("B14") "Barcelona"
Set-TextValue $ws.Range("C14") "Europe"
Set-TextValue $ws.Range("D14") "21st Oct 17:45"
Set-TextValue $ws.Range("G14") "1.18"

# --- Row 15: now Wingate & Finchley vs Ramsgate ---
Set-TextValue $ws.Range("A15") "Wingate & Finchley vs Ramsgate"
Set-TextValue $ws.Range("B15") "Ramsgate"
Set-TextValue $ws.Range("C15") "England"
$ws.Range("D15").ClearContents()
Set-TextValue $ws.Range("E15") "10 of 12 Tips"
Set-TextValue $ws.Range("F15") "83"
Set-TextValue $ws.Range("G15") "2.08"

Write-Host "Applied 2025-10-12 fixtures update"
